$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Agenda")

# Add 0.5 hours of "AI Concepts" (column C) and "Block Chain" (column G) work for "Mon" (row 2)
$ws.Range("C2").Value = 0.5
$ws.Range("G2").Value = 0.5

# Update the selection to match the final state of the sheet
$ws.Activate()
$ws.Range("D21").Select()
